$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 14.82793633333333
$ws.Range("H2").Value = 44.483809
$ws.Range("I2").Value = 0.04809816098739855
$ws.Range("J2").Value = 0.04809816098739855
$ws.Range("M2").Value = 81.76504766666666
$ws.Range("N2").Value = 245.295143
$ws.Range("O2").Value = 0.2487357456271184
$ws.Range("P2").Value = 0.2487357456271185
$ws.Range("Q2").Value = 1212.406921093298
$ws.Range("R2").Value = 10911.66228983969
$ws.Range("S2").Value = 0.01196373193649376
$ws.Range("T2").Value = 0.01196373193649376
$ws.Range("G3").Value = 14.82793633333333
$ws.Range("H3").Value = 44.483809
$ws.Range("I3").Value = 0.04809816098739855
$ws.Range("J3").Value = 0.04809816098739855
$ws.Range("O3").Value = 0.06890783200806287
$ws.Range("P3").Value = 0.06890783200806287
$ws.Range("Q3").Value = 335.8758598748073
$ws.Range("R3").Value = 3022.882738873266
$ws.Range("S3").Value = 0.003314339997216423
$ws.Range("T3").Value = 0.003314339997216423
$ws.Range("G4").Value = 14.82793633333333
$ws.Range("H4").Value = 44.483809
$ws.Range("I4").Value = 0.04809816098739855
$ws.Range("J4").Value = 0.04809816098739855
$ws.Range("M4").Value = 92.136571
$ws.Range("N4").Value = 276.409713
$ws.Range("O4").Value = 0.2802867403763996
$ws.Range("P4").Value = 0.2802867403763997
$ws.Range("Q4").Value = 1366.195208759646
$ws.Range("R4").Value = 12295.75687883682
$ws.Range("S4").Value = 0.01348127676125725
$ws.Range("T4").Value = 0.01348127676125725
$ws.Range("G5").Value = 14.82793633333333
$ws.Range("H5").Value = 44.483809
$ws.Range("I5").Value = 0.04809816098739855
$ws.Range("J5").Value = 0.04809816098739855
$ws.Range("M5").Value = 4.300069
$ws.Range("N5").Value = 12.900207
$ws.Range("O5").Value = 0.01308115019174747
$ws.Range("P5").Value = 0.01308115019174747
$ws.Range("Q5").Value = 63.76114936094033
$ws.Range("R5").Value = 573.850344248463
$ws.Range("S5").Value = 0.0006291792678230091
$ws.Range("T5").Value = 0.0006291792678230092
$ws.Range("G6").Value = 14.82793633333333
$ws.Range("H6").Value = 44.483809
$ws.Range("I6").Value = 0.04809816098739855
$ws.Range("J6").Value = 0.04809816098739855
$ws.Range("M6").Value = 127.8693006666667
$ws.Range("N6").Value = 383.607902
$ws.Range("O6").Value = 0.3889885317966715
$ws.Range("P6").Value = 0.3889885317966715
$ws.Range("Q6").Value = 1896.037849273191
$ws.Range("R6").Value = 17064.34064345872
$ws.Range("S6").Value = 0.0187096330246081
$ws.Range("T6").Value = 0.01870963302460811
$ws.Range("I7").Value = 0.007611361220195336
$ws.Range("J7").Value = 0.007611361220195337
$ws.Range("M7").Value = 81.76504766666666
$ws.Range("N7").Value = 245.295143
$ws.Range("O7").Value = 0.2487357456271184
$ws.Range("P7").Value = 0.2487357456271185
$ws.Range("Q7").Value = 191.8590406132921
$ws.Range("R7").Value = 1726.731365519629
$ws.Range("S7").Value = 0.001893217608342621
$ws.Range("T7").Value = 0.001893217608342621
$ws.Range("I8").Value = 0.007611361220195336
$ws.Range("J8").Value = 0.007611361220195337
$ws.Range("O8").Value = 0.06890783200806287
$ws.Range("P8").Value = 0.06890783200806287
$ws.Range("S8").Value = 0.0005244824003139047
$ws.Range("T8").Value = 0.0005244824003139047
$ws.Range("I9").Value = 0.007611361220195336
$ws.Range("J9").Value = 0.007611361220195337
$ws.Range("M9").Value = 92.136571
$ws.Range("N9").Value = 276.409713
$ws.Range("O9").Value = 0.2802867403763996
$ws.Range("P9").Value = 0.2802867403763997
$ws.Range("Q9").Value = 216.1954847690377
$ws.Range("R9").Value = 1945.759362921339
$ws.Range("S9").Value = 0.002133363626235887
$ws.Range("T9").Value = 0.002133363626235887
$ws.Range("I10").Value = 0.007611361220195336
$ws.Range("J10").Value = 0.007611361220195337
$ws.Range("M10").Value = 4.300069
$ws.Range("N10").Value = 12.900207
$ws.Range("O10").Value = 0.01308115019174747
$ws.Range("P10").Value = 0.01308115019174747
$ws.Range("Q10").Value = 10.08997287293567
$ws.Range("R10").Value = 90.809755856421
$ws.Range("S10").Value = 0.00009956535928501747
$ws.Range("T10").Value = 0.0000995653592850175
$ws.Range("I11").Value = 0.007611361220195336
$ws.Range("J11").Value = 0.007611361220195337
$ws.Range("M11").Value = 127.8693006666667
$ws.Range("N11").Value = 383.607902
$ws.Range("O11").Value = 0.3889885317966715
$ws.Range("P11").Value = 0.3889885317966715
$ws.Range("Q11").Value = 300.0411795736118
$ws.Range("R11").Value = 2700.370616162506
$ws.Range("S11").Value = 0.002960732226017906
$ws.Range("T11").Value = 0.002960732226017906
$ws.Range("G12").Value = 134.1796616666666
$ws.Range("H12").Value = 402.538985
$ws.Range("I12").Value = 0.4352456621741633
$ws.Range("J12").Value = 0.4352456621741634
$ws.Range("M12").Value = 81.76504766666666
$ws.Range("N12").Value = 245.295143
$ws.Range("O12").Value = 0.2487357456271184
$ws.Range("P12").Value = 0.2487357456271185
$ws.Range("Q12").Value = 10971.2064320722
$ws.Range("R12").Value = 98740.85788864984
$ws.Range("S12").Value = 0.1082611543118594
$ws.Range("T12").Value = 0.1082611543118594
$ws.Range("G13").Value = 134.1796616666666
$ws.Range("H13").Value = 402.538985
$ws.Range("I13").Value = 0.4352456621741633
$ws.Range("J13").Value = 0.4352456621741634
$ws.Range("O13").Value = 0.06890783200806287
$ws.Range("P13").Value = 0.06890783200806287
$ws.Range("Q13").Value = 3039.378388662876
$ws.Range("R13").Value = 27354.40549796589
$ws.Range("S13").Value = 0.02999183497133533
$ws.Range("T13").Value = 0.02999183497133533
$ws.Range("G14").Value = 134.1796616666666
$ws.Range("H14").Value = 402.538985
$ws.Range("I14").Value = 0.4352456621741633
$ws.Range("J14").Value = 0.4352456621741634
$ws.Range("M14").Value = 92.136571
$ws.Range("N14").Value = 276.409713
$ws.Range("O14").Value = 0.2802867403763996
$ws.Range("P14").Value = 0.2802867403763997
$ws.Range("Q14").Value = 12362.85392390681
$ws.Range("R14").Value = 111265.6853151613
$ws.Range("S14").Value = 0.1219935879137638
$ws.Range("T14").Value = 0.1219935879137639
$ws.Range("G15").Value = 134.1796616666666
$ws.Range("H15").Value = 402.538985
$ws.Range("I15").Value = 0.4352456621741633
$ws.Range("J15").Value = 0.4352456621741634
$ws.Range("M15").Value = 4.300069
$ws.Range("N15").Value = 12.900207
$ws.Range("O15").Value = 0.01308115019174747
$ws.Range("P15").Value = 0.01308115019174747
$ws.Range("Q15").Value = 576.9818035633216
$ws.Range("R15").Value = 5192.836232069894
$ws.Range("S15").Value = 0.00569351387720681
$ws.Range("T15").Value = 0.005693513877206812
$ws.Range("G16").Value = 134.1796616666666
$ws.Range("H16").Value = 402.538985
$ws.Range("I16").Value = 0.4352456621741633
$ws.Range("J16").Value = 0.4352456621741634
$ws.Range("M16").Value = 127.8693006666667
$ws.Range("N16").Value = 383.607902
$ws.Range("O16").Value = 0.3889885317966715
$ws.Range("P16").Value = 0.3889885317966715
$ws.Range("Q16").Value = 17157.4595010066
$ws.Range("R16").Value = 154417.1355090594
$ws.Range("S16").Value = 0.1693055710999978
$ws.Range("T16").Value = 0.1693055710999979
$ws.Range("G17").Value = 0.961127
$ws.Range("H17").Value = 2.883381
$ws.Range("I17").Value = 0.003117658461441694
$ws.Range("J17").Value = 0.003117658461441695
$ws.Range("M17").Value = 81.76504766666666
$ws.Range("N17").Value = 245.295143
$ws.Range("O17").Value = 0.2487357456271184
$ws.Range("P17").Value = 0.2487357456271185
$ws.Range("Q17").Value = 78.58659496872032
$ws.Range("R17").Value = 707.279354718483
$ws.Range("S17").Value = 0.0007754731020173948
$ws.Range("T17").Value = 0.000775473102017395
$ws.Range("G18").Value = 0.961127
$ws.Range("H18").Value = 2.883381
$ws.Range("I18").Value = 0.003117658461441694
$ws.Range("J18").Value = 0.003117658461441695
$ws.Range("O18").Value = 0.06890783200806287
$ws.Range("P18").Value = 0.06890783200806287
$ws.Range("Q18").Value = 21.771023985866
$ws.Range("R18").Value = 195.939215872794
$ws.Range("S18").Value = 0.00021483108551954
$ws.Range("T18").Value = 0.0002148310855195401
$ws.Range("G19").Value = 0.961127
$ws.Range("H19").Value = 2.883381
$ws.Range("I19").Value = 0.003117658461441694
$ws.Range("J19").Value = 0.003117658461441695
$ws.Range("M19").Value = 92.136571
$ws.Range("N19").Value = 276.409713
$ws.Range("O19").Value = 0.2802867403763996
$ws.Range("P19").Value = 0.2802867403763997
$ws.Range("Q19").Value = 88.554946075517
$ws.Range("R19").Value = 796.994514679653
$ws.Range("S19").Value = 0.0008738383277643937
$ws.Range("T19").Value = 0.000873838327764394
$ws.Range("G20").Value = 0.961127
$ws.Range("H20").Value = 2.883381
$ws.Range("I20").Value = 0.003117658461441694
$ws.Range("J20").Value = 0.003117658461441695
$ws.Range("M20").Value = 4.300069
$ws.Range("N20").Value = 12.900207
$ws.Range("O20").Value = 0.01308115019174747
$ws.Range("P20").Value = 0.01308115019174747
$ws.Range("Q20").Value = 4.132912417762999
$ws.Range("R20").Value = 37.196211759867
$ws.Range("S20").Value = 0.00004078255858069114
$ws.Range("T20").Value = 0.00004078255858069115
$ws.Range("G21").Value = 0.961127
$ws.Range("H21").Value = 2.883381
$ws.Range("I21").Value = 0.003117658461441694
$ws.Range("J21").Value = 0.003117658461441695
$ws.Range("M21").Value = 127.8693006666667
$ws.Range("N21").Value = 383.607902
$ws.Range("O21").Value = 0.3889885317966715
$ws.Range("P21").Value = 0.3889885317966715
$ws.Range("Q21").Value = 122.8986373418513
$ws.Range("R21").Value = 1106.087736076662
$ws.Range("S21").Value = 0.001212733387559674
$ws.Range("T21").Value = 0.001212733387559675
$ws.Range("G22").Value = 155.9696986666667
$ws.Range("H22").Value = 467.909096
$ws.Range("I22").Value = 0.505927157156801
$ws.Range("J22").Value = 0.505927157156801
$ws.Range("M22").Value = 81.76504766666666
$ws.Range("N22").Value = 245.295143
$ws.Range("O22").Value = 0.2487357456271184
$ws.Range("P22").Value = 0.2487357456271185
$ws.Range("Q22").Value = 12752.86984603564
$ws.Range("R22").Value = 114775.8286143207
$ws.Range("S22").Value = 0.1258421686684053
$ws.Range("T22").Value = 0.1258421686684053
$ws.Range("G23").Value = 155.9696986666667
$ws.Range("H23").Value = 467.909096
$ws.Range("I23").Value = 0.505927157156801
$ws.Range("J23").Value = 0.505927157156801
$ws.Range("O23").Value = 0.06890783200806287
$ws.Range("P23").Value = 0.06890783200806287
$ws.Range("Q23").Value = 3532.956675590522
$ws.Range("R23").Value = 31796.6100803147
$ws.Range("S23").Value = 0.03486234355367767
$ws.Range("T23").Value = 0.03486234355367767
$ws.Range("G24").Value = 155.9696986666667
$ws.Range("H24").Value = 467.909096
$ws.Range("I24").Value = 0.505927157156801
$ws.Range("J24").Value = 0.505927157156801
$ws.Range("M24").Value = 92.136571
$ws.Range("N24").Value = 276.409713
$ws.Range("O24").Value = 0.2802867403763996
$ws.Range("P24").Value = 0.2802867403763997
$ws.Range("Q24").Value = 14370.51321504994
$ws.Range("R24").Value = 129334.6189354494
$ws.Range("S24").Value = 0.1418046737473782
$ws.Range("T24").Value = 0.1418046737473782
$ws.Range("G25").Value = 155.9696986666667
$ws.Range("H25").Value = 467.909096
$ws.Range("I25").Value = 0.505927157156801
$ws.Range("J25").Value = 0.505927157156801
$ws.Range("M25").Value = 4.300069
$ws.Range("N25").Value = 12.900207
$ws.Range("O25").Value = 0.01308115019174747
$ws.Range("P25").Value = 0.01308115019174747
$ws.Range("Q25").Value = 670.6804661758746
$ws.Range("R25").Value = 6036.124195582872
$ws.Range("S25").Value = 0.006618109128851939
$ws.Range("T25").Value = 0.00661810912885194
$ws.Range("G26").Value = 155.9696986666667
$ws.Range("H26").Value = 467.909096
$ws.Range("I26").Value = 0.505927157156801
$ws.Range("J26").Value = 0.505927157156801
$ws.Range("M26").Value = 127.8693006666667
$ws.Range("N26").Value = 383.607902
$ws.Range("O26").Value = 0.3889885317966715
$ws.Range("P26").Value = 0.3889885317966715
$ws.Range("Q26").Value = 19943.7362936974
$ws.Range("R26").Value = 179493.6266432766
$ws.Range("S26").Value = 0.1967998620584879
$ws.Range("T26").Value = 0.1967998620584879
